# Update the October 2023 "next statistics" workbook with the new monthly figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (Circulation, ILL Loans, ILL Borrows)
$data = @{
    3  = @(61103, 9149, 10597)
    4  = @(32676, 3718, 3806)
    5  = @(104840, 9635, 9184)
    6  = @(2390, 936, 228)
    7  = @(67576, 11236, 8857)
    8  = @(7175, 1894, 1516)
    9  = @(8237, 1600, 918)
    10 = @(3812, 561, 353)
    11 = @(425, 343, 30)
    12 = @(3, 0, 0)
    13 = @(1414, 340, 385)
    14 = @(4291, 1769, 1353)
    15 = @(6965, 2695, 1113)
    16 = @(5016, 1891, 810)
    17 = @(3057, 1233, 253)
    18 = @(24379, 3796, 4677)
    19 = @(2096, 870, 554)
    20 = @(25865, 3264, 4330)
    21 = @(424, 538, 32)
    22 = @(23901, 3333, 4144)
    23 = @(1560, 627, 250)
    24 = @(27860, 3622, 5395)
    25 = @(109363, 10448, 13201)
    26 = @(8507, 2791, 1365)
    27 = @(0, 0, 0)
    28 = @(7452, 1588, 1750)
    29 = @(1966, 607, 421)
    30 = @(20586, 3728, 3673)
    31 = @(645, 248, 306)
    32 = @(3786, 2398, 573)
    33 = @(22869, 4488, 3976)
    34 = @(14578, 4147, 3080)
    35 = @(7951, 890, 1860)
    36 = @(79846, 7854, 7990)
    37 = @(11702, 3818, 1686)
    38 = @(35372, 2703, 3960)
    39 = @(1402, 1322, 235)
    40 = @(2762, 651, 979)
    41 = @(3987, 453, 161)
    42 = @(14253, 735, 439)
    43 = @(371, 141, 75)
    44 = @(1186, 108, 114)
    45 = @(1045, 14, 7)
    46 = @(4554, 1271, 598)
    47 = @(17289, 4947, 2983)
    48 = @(42906, 4895, 6044)
    49 = @(20712, 4870, 1831)
    50 = @(15686, 1769, 2367)
    51 = @(44407, 4125, 6607)
    52 = @(6768, 930, 1667)
    53 = @(19086, 4126, 3272)
    54 = @(2804, 1854, 1047)
    55 = @(3104, 1764, 201)
    56 = @(5541, 1410, 1766)
    57 = @(17686, 6695, 3661)
    58 = @(20842, 1451, 796)
    59 = @(947551, 144143, 133019)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}

# Update the active cell/selection on the frozen (bottom-right) pane to A2.
$ws.Range("A2").Select()
